$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the "Электро" (row 4) and "Психиатр" (row 5) line items entirely.
# Deleting the whole rows shifts everything below up by two rows and drops
# the now-unused shared strings from the shared string table.
$ws.Rows("4:5").Delete()

# Update the "Связь" line (row 3) amounts.
$ws.Range("B3").Value = 140
$ws.Range("C3").Value = 0

# Add a (currently empty) statistics column D next to the existing data,
# giving it its own formatting so the sheet's used range grows to D8.
$ws.Range("D3:D8").HorizontalAlignment = 1

# Leave the selection where the author left it when they saved.
[void]$ws.Range("G15").Select()
